# The presentation ships two theme parts:
#   ppt/theme/theme1.xml  -> used by the slide master (was the "Integral" theme)
#   ppt/theme/theme2.xml  -> used by the notes master (was the "Office Theme")
#
# The commit swaps the content of the two theme parts: the slide master's
# theme becomes the stock "Office Theme" colour scheme, and the notes
# master's theme becomes the old "Integral" colour scheme.
#
# The PowerPoint object model reaches the slide-facing theme (theme1.xml)
# through ThemeColorScheme on a Slide (any slide works - they all share the
# single presentation theme). Each of the 12 slots maps 1:1 onto the OOXML
# <a:clrScheme> children in the standard order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
#
# RGB is written as a standard VBA RGB() integer (R + G*256 + B*65536), i.e.
# the low byte is Red, not Blue.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function Set-ThemeColor($index, $r, $g, $b) {
    $tcs.Item($index).RGB = $r + ($g * 256) + ($b * 65536)
}

# New "Office Theme" colour scheme for the slide master's theme part.
Set-ThemeColor 1  0x00 0x00 0x00   # dk1
Set-ThemeColor 2  0xFF 0xFF 0xFF   # lt1
Set-ThemeColor 3  0x44 0x54 0x6A   # dk2
Set-ThemeColor 4  0xE7 0xE6 0xE6   # lt2
Set-ThemeColor 5  0x5B 0x9B 0xD5   # accent1
Set-ThemeColor 6  0xED 0x7D 0x31   # accent2
Set-ThemeColor 7  0xA5 0xA5 0xA5   # accent3
Set-ThemeColor 8  0xFF 0xC0 0x00   # accent4
Set-ThemeColor 9  0x44 0x72 0xC4   # accent5
Set-ThemeColor 10 0x70 0xAD 0x47   # accent6
Set-ThemeColor 11 0x05 0x63 0xC1   # hlink
Set-ThemeColor 12 0x95 0x4F 0x72   # folHlink
